# Add files via upload
# Rebuild the player table (A2:C19) with the updated roster / positions / teams.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Ja Morant",        "PG",          "Memphis Grizzlies"),
    @("Scottie Barnes",   "PG,SG,SF,PF", "Toronto Raptors"),
    @("De'Aaron Fox",     "PG,SG",       "San Antonio Spurs"),
    @("Scoot Henderson",  "PG",          "Portland Trail Blazers"),
    @("DeMar DeRozan",    "SF,PF",       "Sacramento Kings"),
    @("Miles Bridges",    "SF,PF",       "Charlotte Hornets"),
    @("Christian Braun",  "SG,SF",       "Denver Nuggets"),
    @("Josh Giddey",      "PG,SG,SF",    "Chicago Bulls"),
    @("Moussa Diabate",   "C",           "Charlotte Hornets"),
    @("Zach Collins",     "PF,C",        "Chicago Bulls"),
    @("Luka Doncic",      "PG,SG",       "Los Angeles Lakers"),
    @("Mikal Bridges",    "SG,SF,PF",    "New York Knicks"),
    @("Isaiah Collier",   "PG,SG",       "Utah Jazz"),
    @("Brook Lopez",      "C",           "Milwaukee Bucks"),
    @("Evan Mobley",      "PF,C",        "Cleveland Cavaliers"),
    @("Nikola Vucevic",   "PF,C",        "Chicago Bulls"),
    @("P.J. Washington",  "SF,PF",       "Dallas Mavericks"),
    @("Tyler Herro",      "PG,SG",       "Miami Heat")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
